$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-Text "2025-10-27 Monday" "2025-10-28 Tuesday"

Replace-Text "852×2=" "402×2="
Replace-Text "544×3=" "129×4="
Replace-Text "225×2=" "970×7="
Replace-Text "622×8=" "874×7="
Replace-Text "838×8=" "773×7="

Replace-Text "219×9=" "168×2="
Replace-Text "448×6=" "443×9="
Replace-Text "364×7=" "828×4="
Replace-Text "399×7=" "757×5="
Replace-Text "904×2=" "943×3="

Replace-Text "669×6=" "947×4="
Replace-Text "128×4=" "465×5="
Replace-Text "169×7=" "959×9="
Replace-Text "143×4=" "693×3="
Replace-Text "620×7=" "772×6="

Replace-Text "619×6=" "796×8="
Replace-Text "305×7=" "324×3="
Replace-Text "321×9=" "902×6="
Replace-Text "857×9=" "871×5="
Replace-Text "516×2=" "188×8="

Replace-Text "522×7=" "551×5="
Replace-Text "265×8=" "560×6="
Replace-Text "182×6=" "376×4="
Replace-Text "182×2=" "545×5="
Replace-Text "878×9=" "354×7="
